# Regenerate the "K" column (G) values for the save_data sheet.
# These values replace the previous "Strike#" derived numbers with
# actual strikeout counts (K), as part of regenerating std/mean and
# writing s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 4
    3  = 4
    4  = 4
    5  = 3
    6  = 4
    7  = 3
    8  = 3
    9  = 7
    10 = 3
    11 = 3
    12 = 3
    13 = 3
    14 = 2
    15 = 3
    16 = 4
    17 = 1
    18 = 3
    19 = 5
    20 = 6
    21 = 3
    22 = 5
    23 = 6
    24 = 4
    25 = 3
    26 = 3
    27 = 4
    28 = 5
    29 = 5
    30 = 3
    31 = 3
    32 = 4
    33 = 5
    34 = 3
    35 = 2
    36 = 5
    37 = 4
    38 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
